$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 4) down onto the
# new row 5 so the new row picks up the same cell style (fill, etc.)
# that the other data rows use, instead of Excel inventing a new style.
$ws.Range("A4:V4").Copy()
$ws.Range("A5:V5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row's contents.
$ws.Range("A5").Value2 = "GMHO:0000191"
$ws.Range("B5").Value2 = "repeated measures study design"
$ws.Range("C5").Value2 = "A study design in which measurements taken on the same study participants at two or more different times in different circumstances  are compared."
$ws.Range("D5").Value2 = "study design"
$ws.Range("S5").Value2 = "Proposed"
$ws.Range("V5").Value2 = "BG"
